# Apply the diff described for the "artfynd" workbook:
#  1. Bump the "Taxonsorteringsordning" (column B) value 79243 -> 79244 on
#     every row that currently holds 79243.
#  2. Bump column B value 91829 -> 91830 on the one row that holds it.
#  3. Swap the full record content of row 12 (Id 130981909 / Tjäder) and
#     row 13 (Id 130981935 / Garnlav) - the source list got re-sorted so
#     the two records exchanged rows. The Garnlav record additionally
#     receives the 79243 -> 79244 bump described in step 1 (its
#     Taxonsorteringsordning was 79243 before the edit); the Tjäder
#     record's B value (57073) is untouched by that rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1/2: simple value bumps on column B ("Taxonsorteringsordning")
# ---------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row   # xlUp = -4162
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq 79243) {
        $cell.Value = 79244
    } elseif ($val -eq 91829) {
        $cell.Value = 91830
    }
}

# ---------------------------------------------------------------------
# Step 3: swap the two records currently sitting in row 12 and row 13.
# Every populated column of each row is captured first, then the
# captured values are written crosswise (row13 -> row12, row12 -> row13).
# Columns Y/AA hold dates stored as plain text (e.g. "2026-01-31"); they
# are written with the cell pre-formatted as text so Excel does not
# reinterpret them as real date serial numbers, and the style is reset
# back to Normal afterwards so no stray formatting is left behind.
# ---------------------------------------------------------------------
$cols = @(1,2,4,5,6,7,8,9,11,12,13,14,16,17,18,19,20,21,22,23,25,27,30,31,33,46,49,50,51)
# A=1 B=2 D=4 E=5 F=6 G=7 H=8 I=9 K=11 L=12 M=13 N=14 P=16 Q=17 R=18 S=19
# T=20 U=21 V=22 W=23 Y=25 AA=27 AD=30 AE=31 AG=33 AT=46 AW=49 AX=50 AY=51
$dateCols = @(25, 27)   # Y, AA

$row12 = @{}
$row13 = @{}
foreach ($c in $cols) {
    $row12[$c] = $ws.Cells.Item(12, $c).Value2
    $row13[$c] = $ws.Cells.Item(13, $c).Value2
}

foreach ($c in $cols) {
    $dest12 = $ws.Cells.Item(12, $c)
    $dest13 = $ws.Cells.Item(13, $c)

    if ($dateCols -contains $c) {
        $dest12.NumberFormat = "@"
        $dest13.NumberFormat = "@"
    }

    $dest12.Value = $row13[$c]
    $dest13.Value = $row12[$c]

    if ($dateCols -contains $c) {
        $dest12.Style = "Normal"
        $dest13.Style = "Normal"
    }
}

# The Garnlav record (now sitting on row 12) still needs the
# 79243 -> 79244 bump from step 1, since it was on row 13 (already
# processed) at the time the loop above ran.
if ($ws.Cells.Item(12, 2).Value2 -eq 79243) {
    $ws.Cells.Item(12, 2).Value = 79244
}
